$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 47.72
$ws.Range("E3").Value = 1.69

$ws.Range("D4").Value = 1.1

$ws.Range("D5").Value = 2.3
$ws.Range("E5").Value = 0.52

$ws.Range("D6").Value = 3.3

$ws.Range("D9").Value = 1.31

$ws.Range("D10").Value = 3.6

$ws.Range("D11").Value = 15.87
$ws.Range("E11").Value = 1.23

$ws.Range("D14").Value = 0.47
$ws.Range("E14").Value = 0.17

$ws.Range("D18").Value = 1.71

$ws.Range("D22").Value = 0.87

$ws.Range("D24").Value = 2.37
$ws.Range("E24").Value = 0.53

$ws.Range("D26").Value = 0.74
